# Update "想去人数" (column F) counts on the 展览 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 30
$ws1.Range("F5").Value = 3324
$ws1.Range("F6").Value = 2126
$ws1.Range("F7").Value = 405
$ws1.Range("F9").Value = 36
$ws1.Range("F10").Value = 20
$ws1.Range("F11").Value = 1228
$ws1.Range("F13").Value = 1311
$ws1.Range("F14").Value = 105

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 30
$ws4.Range("F5").Value = 3324
$ws4.Range("F6").Value = 2126
$ws4.Range("F7").Value = 405
$ws4.Range("F10").Value = 36
$ws4.Range("F11").Value = 20
$ws4.Range("F14").Value = 1228
$ws4.Range("F16").Value = 1311
$ws4.Range("F17").Value = 105
